# Adds the latest daily COVID positive-test-by-district rows
# (date serial 44195 = 2020-12-30) to the bottom of the data table,
# matching the "Updated: st 31. 12. 2020" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append: @(DateSerial, District, DailyPositiveTests)
$newData = @(
    @(44195, 'Bánovce nad Bebravou', 21),
    @(44195, 'Banská Bystrica', 184),
    @(44195, 'Banská Štiavnica', 5),
    @(44195, 'Bardejov', 79),
    @(44195, 'Bratislava', 487),
    @(44195, 'Brezno', 44),
    @(44195, 'Bytča', 58),
    @(44195, 'Čadca', 16),
    @(44195, 'Detva', 36),
    @(44195, 'Dolný Kubín', 27),
    @(44195, 'Dunajská Streda', 119),
    @(44195, 'Galanta', 73),
    @(44195, 'Gelnica', 14),
    @(44195, 'Hlohovec', 40),
    @(44195, 'Humenné', 52),
    @(44195, 'Ilava', 118),
    @(44195, 'Kežmarok', 34),
    @(44195, 'Komárno', 33),
    @(44195, 'Košice', 268),
    @(44195, 'Košice - okolie', 108),
    @(44195, 'Krupina', 18),
    @(44195, 'Kysucké Nové Mesto', 34),
    @(44195, 'Levice', 82),
    @(44195, 'Levoča', 28),
    @(44195, 'Liptovský Mikuláš', 112),
    @(44195, 'Lučenec', 75),
    @(44195, 'Malacky', 73),
    @(44195, 'Martin', 260),
    @(44195, 'Medzilaborce', 6),
    @(44195, 'Michalovce', 53),
    @(44195, 'Myjava', 56),
    @(44195, 'Námestovo', 13),
    @(44195, 'Nitra', 584),
    @(44195, 'Nové Mesto nad Váhom', 80),
    @(44195, 'Nové Zámky', 105),
    @(44195, 'Partizánske', 32),
    @(44195, 'Pezinok', 33),
    @(44195, 'Piešťany', 91),
    @(44195, 'Poltár', 21),
    @(44195, 'Poprad', 96),
    @(44195, 'Považská Bystrica', 150),
    @(44195, 'Prešov', 143),
    @(44195, 'Prievidza', 102),
    @(44195, 'Púchov', 94),
    @(44195, 'Revúca', 18),
    @(44195, 'Rimavská Sobota', 74),
    @(44195, 'Rožňava', 29),
    @(44195, 'Ružomberok', 107),
    @(44195, 'Sabinov', 42),
    @(44195, 'Senec', 57),
    @(44195, 'Senica', 150),
    @(44195, 'Skalica', 165),
    @(44195, 'Snina', 27),
    @(44195, 'Sobrance', 14),
    @(44195, 'Spišská Nová Ves', 37),
    @(44195, 'Stará Ľubovňa', 62),
    @(44195, 'Stropkov', 16),
    @(44195, 'Svidník', 66),
    @(44195, 'Šaľa', 41),
    @(44195, 'Topoľčany', 180),
    @(44195, 'Trebišov', 67),
    @(44195, 'Trenčín', 183),
    @(44195, 'Trnava', 147),
    @(44195, 'Turčianske Teplice', 23),
    @(44195, 'Tvrdošín', 11),
    @(44195, 'Veľký Krtíš', 19),
    @(44195, 'Vranov nad Topľou', 78),
    @(44195, 'Zlaté Moravce', 69),
    @(44195, 'Zvolen', 65),
    @(44195, 'Žarnovica', 20),
    @(44195, 'Žiar nad Hronom', 14),
    @(44195, 'Žilina', 377)
)

# Find the first empty row right after the existing data block in column A
# (xlUp = -4162), so the script appends in place regardless of current extent.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $entry = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
}

$newLastRow = $startRow + $newData.Count - 1
Write-Output "Appended $($newData.Count) rows ($startRow..$newLastRow)"